$wb = $excel.ActiveWorkbook

# Use "Greece" worksheet as the structural template: it already has the
# exact row heights / column widths / styles used by the new country
# sheets being added (Norway, Poland).
$template = $wb.Worksheets.Item("Greece")

# --- Add "Norway" sheet (copied after the last sheet, i.e. after "Hungary") ---
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Activate()
$norway.Range("H26").Select()
$norway.Range("B4").Value = "NGC-2931/T3058"
$norway.Range("B2").Value = "Norway Market"

# --- Add "Poland" sheet (copied after "Norway") ---
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$poland = $wb.Worksheets.Item($wb.Worksheets.Count)
$poland.Name = "Poland"
$poland.Activate()
$poland.Range("H26").Select()
$poland.Range("B4").Value = "NGC-2920/T3101"
$poland.Range("B2").Value = "Poland Market"

# "Norway" ends up being the active tab in the final workbook.
$norway.Activate()
